$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values pulled from the refreshed scheduled-operations query.
# These columns are stored as text in the report, so assign string values
# (quoted) to keep them as shared-string text cells rather than numbers.
$ws.Range("D3").Value  = "3049"
$ws.Range("G3").Value  = "6023"

$ws.Range("D4").Value  = "7524"
$ws.Range("G4").Value  = "8245"

$ws.Range("D8").Value  = "346"
$ws.Range("G8").Value  = "474"

$ws.Range("D9").Value  = "5200"
$ws.Range("G9").Value  = "6737"

$ws.Range("D10").Value = "598"
$ws.Range("G10").Value = "825"

$ws.Range("C12").Value = "970"
$ws.Range("D12").Value = "719"
$ws.Range("G12").Value = "1763"

$ws.Range("C15").Value = "20003"
$ws.Range("D15").Value = "105322"
$ws.Range("G15").Value = "155850"

$ws.Range("C18").Value = "1825"
$ws.Range("D18").Value = "30003"
$ws.Range("G18").Value = "43838"

$ws.Range("D20").Value = "1599"
$ws.Range("G20").Value = "2872"

$ws.Range("D21").Value = "1581"
$ws.Range("G21").Value = "2467"

$ws.Range("C22").Value = "1934"
$ws.Range("G22").Value = "22319"

$ws.Range("C27").Value = "31201"
$ws.Range("D27").Value = "196256"
$ws.Range("G27").Value = "302415"

# Refresh the "best fit" width on column D (CC) to match the new data.
$ws.Columns.Item(4).AutoFit() | Out-Null
